$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.208.48'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '3.137.12'
$ws.Range('E3').Value = '  +3.50%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '580.31'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').Value = '174.85'
$ws.Range('E6').Value = '  +3.86%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.131.88'
$ws.Range('E8').Value = '  +3.39%  '
$ws.Range('D9').Value = '0.524'
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').Value = '6.51'
$ws.Range('E10').Value = '  -2.57%  '
$ws.Range('E11').Value = '  +2.11%  '
$ws.Range('E12').Value = '  -1.04%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000250'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('E14').Value = '  +1.94%  '
$ws.Range('E15').Value = '  -0.74%  '
$ws.Range('D16').Value = '3.656.84'
$ws.Range('E16').Value = '  +3.51%  '
$ws.Range('D17').Value = '67.267.47'
$ws.Range('E17').Value = '  +1.45%  '
$ws.Range('E18').Value = '  -0.37%  '
$ws.Range('D19').Value = '3.136.97'
$ws.Range('E19').Value = '  +3.47%  '
$ws.Range('D20').Value = '16.17'
$ws.Range('D21').Value = '488.36'
$ws.Range('E21').Value = '  +4.59%  '
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('D23').Value = '7.68'
$ws.Range('E23').Value = '  +3.87%  '
$ws.Range('D24').Value = '84.27'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('D25').Value = '13.27'
$ws.Range('E25').Value = '  +4.25%  '
$ws.Range('D26').Value = '2.33'
$ws.Range('E26').Value = '  +3.30%  '
$ws.Range('D27').Value = '10.07'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('E29').Value = '  -2.74%  '
$ws.Range('D30').Value = '2.41'
$ws.Range('E30').Value = '  -0.51%  '
$ws.Range('D31').Value = '2.69'
$ws.Range('E31').Value = '  +1.53%  '
$ws.Range('D32').Value = '28.85'
$ws.Range('E32').Value = '  +2.43%  '
$ws.Range('D33').Value = '0.0₃0996'
$ws.Range('E33').Value = '  -0.03%  '
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  +1.53%  '
$ws.Range('D37').Value = '0.991'
$ws.Range('E37').Value = '  +0.06%  '
$ws.Range('D38').Value = '47.73'
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('D39').Value = '2.11'
$ws.Range('E39').Value = '  +2.68%  '
$ws.Range('D40').Value = '50.14'
$ws.Range('E40').Value = '  +1.27%  '
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('E42').Value = '  +1.55%  '
$ws.Range('D43').Value = '8.67'
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').Value = '2.79'
$ws.Range('E44').Value = '  -1.72%  '
$ws.Range('D45').Value = '2.844.54'
$ws.Range('E45').Value = '  +5.23%  '
$ws.Range('E46').Value = '  -0.33%  '
$ws.Range('D47').Value = '384.18'
$ws.Range('E47').Value = '  +1.39%  '
$ws.Range('D48').Value = '135.56'
$ws.Range('E48').Value = '  +0.99%  '
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').Value = '24.94'
$ws.Range('E50').Value = '  +2.09%  '
$ws.Range('D51').Value = '2.22'
$ws.Range('E51').Value = '  -0.47%  '
